$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder / update the bottom three rows of the ranking table.
# Row 41: XGBoost: C Si   -> 96.8  (unchanged value, label moved up from row 42)
# Row 42: RFC: C Si N     -> 96.9  (was row 43's label, value bumped from 96.8)
# Row 43: RFC: C Si       -> 97.3  (was row 41's label, value bumped from 96.9)
$ws.Range("A41").Value = "XGBoost: C Si"
$ws.Range("B41").Value = 96.8

$ws.Range("A42").Value = "RFC: C Si N"
$ws.Range("B42").Value = 96.9

$ws.Range("A43").Value = "RFC: C Si"
$ws.Range("B43").Value = 97.3

# Restore the cursor/selection and scroll position to match the saved view.
$ws.Range("B44").Select()
